$wb = $excel.ActiveWorkbook
$aw = $excel.ActiveWindow
$aw.Left = -110
$aw.Top = -110
$aw.Width = 18020
$aw.Height = 11020
